$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 1206
$ws.Range("J103").Value = 1299.8
$ws.Range("L103").Value = 3899.4
$ws.Range("N103").Value = -5071.4

# Row 121
$ws.Range("H121").Value = 1362.5
$ws.Range("I121").Value = 1525
$ws.Range("K121").Value = 4575
$ws.Range("M121").Value = -2828

# Row 129
$ws.Range("H129").Value = 808.8421
$ws.Range("I129").Value = 333.66666
$ws.Range("J129").Value = 935.55554
$ws.Range("K129").Value = 1000.99998
$ws.Range("L129").Value = 2806.66662
$ws.Range("M129").Value = 3999.00002
$ws.Range("N129").Value = -12806.66662

# Row 137
$ws.Range("H137").Value = 3657.394
$ws.Range("I137").Value = 2138.3333
$ws.Range("K137").Value = 6414.999899999999
$ws.Range("M137").Value = -3864.999899999999

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1965.7
$ws.Range("I45").Value = 1826.6666
$ws.Range("J45").Value = 2382.8
$ws.Range("K45").Value = 1826.6666
$ws.Range("L45").Value = 2382.8
$ws.Range("M45").Value = -1449.6666
$ws.Range("N45").Value = -3136.8

# Row 122
$ws.Range("H122").Value = 1451.2069
$ws.Range("I122").Value = 1211.1904
$ws.Range("J122").Value = 2081.25
$ws.Range("K122").Value = 3633.5712
$ws.Range("L122").Value = 6243.75
$ws.Range("M122").Value = -1183.5712
$ws.Range("N122").Value = -11143.75

$ws = $wb.Worksheets.Item("BSM")
# Row 95
$ws.Range("H95").Value = 29500
$ws.Range("J95").Value = 29500
$ws.Range("L95").Value = 29500
$ws.Range("N95").Value = -34992

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2757780.2
$ws.Range("I58").Value = 4330940.5
$ws.Range("J58").Value = 4750
$ws.Range("K58").Value = 4330940.5
$ws.Range("L58").Value = 4750
$ws.Range("M58").Value = -4330737.5
$ws.Range("N58").Value = -5156

# Row 62
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("M62").Value = -2376

# Row 65
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("N65").ClearContents()

# Row 86
$ws.Range("H86").Value = 3856.8
$ws.Range("I86").Value = 3834.5
$ws.Range("J86").Value = 3914.1428
$ws.Range("K86").Value = 3834.5
$ws.Range("L86").Value = 3914.1428
$ws.Range("M86").Value = -2711.5
$ws.Range("N86").Value = -6160.1428

# Row 89
$ws.Range("H89").Value = 3856.8
$ws.Range("I89").Value = 3834.5
$ws.Range("J89").Value = 3914.1428
$ws.Range("K89").Value = 19172.5
$ws.Range("L89").Value = 19570.714
$ws.Range("M89").Value = -13556.5
$ws.Range("N89").Value = -30802.714

# Row 107
$ws.Range("H107").Value = 562.7778
$ws.Range("I107").Value = 392.90475
$ws.Range("J107").Value = 800.6
$ws.Range("K107").Value = 392.90475
$ws.Range("L107").Value = 800.6
$ws.Range("M107").Value = 1527.09525
$ws.Range("N107").Value = -4640.6

# Row 134
$ws.Range("H134").Value = 2878.2554
$ws.Range("I134").Value = 2235.1538
$ws.Range("J134").Value = 3674.476
$ws.Range("K134").Value = 6705.4614
$ws.Range("L134").Value = 11023.428
$ws.Range("M134").Value = -4170.4614
$ws.Range("N134").Value = -16093.428

# Row 136
$ws.Range("H136").Value = 2757780.2
$ws.Range("I136").Value = 4330940.5
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 12992821.5
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -12990271.5
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("CUL")
# Row 36
$ws.Range("H36").Value = 689.8
$ws.Range("I36").Value = 633
$ws.Range("J36").Value = 775
$ws.Range("K36").Value = 1899
$ws.Range("L36").Value = 2325
$ws.Range("M36").Value = -1730
$ws.Range("N36").Value = -2663

# Row 61
$ws.Range("H61").Value = 479.26315
$ws.Range("I61").Value = 50
$ws.Range("J61").Value = 559.75
$ws.Range("K61").Value = 150
$ws.Range("L61").Value = 1679.25
$ws.Range("M61").Value = 65
$ws.Range("N61").Value = -2109.25

# Row 63
$ws.Range("H63").Value = 2626.2632
$ws.Range("I63").Value = 1749.5
$ws.Range("J63").Value = 2860.0667
$ws.Range("K63").Value = 5248.5
$ws.Range("L63").Value = 8580.2001
$ws.Range("M63").Value = -4499.5
$ws.Range("N63").Value = -10078.2001

# Row 66
$ws.Range("H66").Value = 2626.2632
$ws.Range("I66").Value = 1749.5
$ws.Range("J66").Value = 2860.0667
$ws.Range("K66").Value = 15745.5
$ws.Range("L66").Value = 25740.6003
$ws.Range("M66").Value = -12001.5
$ws.Range("N66").Value = -33228.6003

# Row 69
$ws.Range("H69").Value = 1369.6296
$ws.Range("I69").Value = 450
$ws.Range("J69").Value = 1443.2
$ws.Range("K69").Value = 1350
$ws.Range("L69").Value = 4329.6
$ws.Range("M69").Value = -539
$ws.Range("N69").Value = -5951.6

# Row 72
$ws.Range("H72").Value = 1369.6296
$ws.Range("I72").Value = 450
$ws.Range("J72").Value = 1443.2
$ws.Range("K72").Value = 4050
$ws.Range("L72").Value = 12988.8
$ws.Range("M72").Value = 6
$ws.Range("N72").Value = -21100.8

# Row 80
$ws.Range("H80").Value = 2530
$ws.Range("J80").Value = 2628.5715
$ws.Range("L80").Value = 7885.7145
$ws.Range("N80").Value = -9757.7145

# Row 83
$ws.Range("H83").Value = 2530
$ws.Range("J83").Value = 2628.5715
$ws.Range("L83").Value = 23657.1435
$ws.Range("N83").Value = -33017.1435

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9375.75
$ws.Range("I80").Value = 26000
$ws.Range("J80").Value = 3834.3333
$ws.Range("K80").Value = 26000
$ws.Range("L80").Value = 3834.3333
$ws.Range("M80").Value = -25002
$ws.Range("N80").Value = -5830.3333

# Row 83
$ws.Range("H83").Value = 9375.75
$ws.Range("I83").Value = 26000
$ws.Range("J83").Value = 3834.3333
$ws.Range("K83").Value = 130000
$ws.Range("L83").Value = 19171.6665
$ws.Range("M83").Value = -125008
$ws.Range("N83").Value = -29155.6665

# Row 122
$ws.Range("H122").Value = 3501.111
$ws.Range("I122").Value = 3616.7827
$ws.Range("J122").Value = 2836
$ws.Range("K122").Value = 10850.3481
$ws.Range("L122").Value = 8508
$ws.Range("M122").Value = -8400.348100000001
$ws.Range("N122").Value = -13408

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3946.2222
$ws.Range("I7").Value = 2823.2354
$ws.Range("J7").Value = 5855.3
$ws.Range("K7").Value = 2823.2354
$ws.Range("L7").Value = 5855.3
$ws.Range("M7").Value = -2711.2354
$ws.Range("N7").Value = -6079.3

# Row 22
$ws.Range("H22").Value = 586
$ws.Range("J22").Value = 586
$ws.Range("L22").Value = 586
$ws.Range("N22").Value = -1176

# Row 27
$ws.Range("H27").Value = 586
$ws.Range("J27").Value = 586
$ws.Range("L27").Value = 586
$ws.Range("N27").Value = -800

# Row 40
$ws.Range("H40").Value = 2727.081
$ws.Range("I40").Value = 2628.4
$ws.Range("J40").Value = 3150
$ws.Range("K40").Value = 2628.4
$ws.Range("L40").Value = 3150
$ws.Range("M40").Value = -2492.4
$ws.Range("N40").Value = -3422

# Row 126
$ws.Range("H126").Value = 3946.2222
$ws.Range("I126").Value = 2823.2354
$ws.Range("J126").Value = 5855.3
$ws.Range("K126").Value = 8469.706200000001
$ws.Range("L126").Value = 17565.9
$ws.Range("M126").Value = -5999.706200000001
$ws.Range("N126").Value = -22505.9

# Row 132
$ws.Range("H132").Value = 3537.279
$ws.Range("I132").Value = 3032.5293
$ws.Range("J132").Value = 5444.1113
$ws.Range("K132").Value = 9097.5879
$ws.Range("L132").Value = 16332.3339
$ws.Range("M132").Value = -6567.5879
$ws.Range("N132").Value = -21392.3339

# Row 136
$ws.Range("H136").Value = 5375.1836
$ws.Range("I136").Value = 4110.7036
$ws.Range("J136").Value = 6927.0454
$ws.Range("K136").Value = 12332.1108
$ws.Range("L136").Value = 20781.1362
$ws.Range("M136").Value = -9782.110799999999
$ws.Range("N136").Value = -25881.1362

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3375.5
$ws.Range("I122").Value = 2662.1177
$ws.Range("J122").Value = 5801
$ws.Range("K122").Value = 7986.353099999999
$ws.Range("L122").Value = 17403
$ws.Range("M122").Value = -5536.353099999999
$ws.Range("N122").Value = -22303

# Row 132
$ws.Range("H132").Value = 1838.5172
$ws.Range("I132").Value = 806.7857
$ws.Range("J132").Value = 2801.4666
$ws.Range("K132").Value = 2420.3571
$ws.Range("L132").Value = 8404.399800000001
$ws.Range("M132").Value = 109.6428999999998
$ws.Range("N132").Value = -13464.3998

# Row 136
$ws.Range("H136").Value = 3190.1956
$ws.Range("I136").Value = 2416.6667
$ws.Range("J136").Value = 4640.5625
$ws.Range("K136").Value = 7250.000100000001
$ws.Range("L136").Value = 13921.6875
$ws.Range("M136").Value = -4700.000100000001
$ws.Range("N136").Value = -19021.6875
